# The deck's "date" placeholders (on the slide master and every slide
# layout) cache the text of a datetimeFigureOut field. The author's
# change simply refreshes that cached text from 8/31/2018 to 9/10/2018
# everywhere it appears (slide master + all 11 layouts).

$p = $ppt.ActivePresentation
$oldDate = "8/31/2018"
$newDate = "9/10/2018"
$ppPlaceholderDate = 16

function Update-DateShape($shapes) {
    for ($i = 1; $i -le $shapes.Count; $i++) {
        $shp = $shapes.Item($i)
        $isDate = $false
        try {
            if ($shp.PlaceholderFormat.Type -eq $ppPlaceholderDate) {
                $isDate = $true
            }
        } catch {
        }
        if ($isDate -and $shp.HasTextFrame) {
            $tr = $shp.TextFrame.TextRange
            if ($tr.Text -eq $oldDate -or $tr.Text -eq $newDate) {
                $tr.Text = $newDate
            }
        }
    }
}

# Slide master's own "Date Placeholder" shape.
Update-DateShape $p.SlideMaster.Shapes

# Every slide layout ("custom layout") has its own "Date Placeholder" shape.
$layouts = $p.SlideMaster.CustomLayouts
for ($li = 1; $li -le $layouts.Count; $li++) {
    $layout = $layouts.Item($li)
    Update-DateShape $layout.Shapes
}
